$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.012.56'
$ws.Range('E2').Value = '  -2.04%  '
$ws.Range('D3').Value = '2.351.11'
$ws.Range('E3').Value = '  -5.05%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '472.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('D9').Value = '2.350.86'
$ws.Range('E9').Value = '  -5.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0961'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.40'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.31%  '
$ws.Range('E12').Value = '  -3.98%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '2.757.87'
$ws.Range('E14').Value = '  -5.08%  '
$ws.Range('D15').Value = '55.040.93'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').Value = '2.350.89'
$ws.Range('E18').Value = '  -5.42%  '
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '311.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.56'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E23').Value = '  -3.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '55.89'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('E26').Value = '  -4.70%  '
$ws.Range('E27').Value = '  -5.28%  '
$ws.Range('D28').Value = '2.451.21'
$ws.Range('E28').Value = '  -5.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  -4.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('E34').Value = '  -2.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.01'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.73%  '
$ws.Range('E36').Value = '  -5.31%  '
$ws.Range('E37').Value = '  -4.69%  '
$ws.Range('E38').Value = '  -4.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.48'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.81%  '
$ws.Range('E40').Value = '  +0.54%  '
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('E42').Value = '  -4.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0943'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0525'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.574'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '254.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.33%  '
$ws.Range('E48').Value = '  -3.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.42'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.73'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.90%  '
$ws.Range('D51').Value = '1.777.83'
$ws.Range('E51').Value = '  -6.16%  '
